# Apply cryptos list update (GitHub Actions refresh) to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.525.90'
$ws.Range("E2").Value = '  +3.08%  '

# Row 3
$ws.Range("D3").Value = '3.403.92'
$ws.Range("E3").Value = '  +2.22%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.34'
$ws.Range("E5").Value = '  +2.23%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.95'
$ws.Range("E6").Value = '  +2.26%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  +2.73%  '

# Row 8
$ws.Range("D8").Value = '3.394.96'
$ws.Range("E8").Value = '  +2.24%  '

# Row 9
$ws.Range("E9").Value = '  +0.03%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.174'
$ws.Range("E10").Value = '  +15.03%  '

# Row 11
$ws.Range("E11").Value = '  +3.12%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.94'
$ws.Range("E12").Value = '  +3.01%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000283'
$ws.Range("E13").Value = '  +7.01%  '

# Row 14
$ws.Range("E14").Value = '  +3.03%  '

# Row 15
$ws.Range("D15").Value = '3.953.54'
$ws.Range("E15").Value = '  +2.61%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.33'
$ws.Range("E16").Value = '  +3.20%  '

# Row 17
$ws.Range("D17").Value = '3.406.53'
$ws.Range("E17").Value = '  +2.58%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.118'
$ws.Range("E18").Value = '  +1.66%  '

# Row 19
$ws.Range("D19").Value = '65.552.89'
$ws.Range("E19").Value = '  +3.22%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.89'
$ws.Range("E20").Value = '  +2.51%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.994'
$ws.Range("E21").Value = '  +2.73%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '472.26'
$ws.Range("E22").Value = '  +15.93%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.11'

# Row 24
$ws.Range("E24").Value = '  +2.26%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.57'
$ws.Range("E25").Value = '  +4.54%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.56'
$ws.Range("E26").Value = '  +2.19%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.89'
$ws.Range("E27").Value = '  +3.18%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.88'
$ws.Range("E28").Value = '  +6.09%  '

# Row 29
$ws.Range("E29").Value = '  +3.67%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.03'
$ws.Range("E30").Value = '  +6.92%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.68'
$ws.Range("E31").Value = '  +4.68%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.53'
$ws.Range("E32").Value = '  +2.21%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '581.26'
$ws.Range("E33").Value = '  +1.68%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '62.30'
$ws.Range("E34").Value = '  +8.92%  '

# Row 35
$ws.Range("E35").Value = '  +2.57%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.08%  '

# Row 37
$ws.Range("E37").Value = '  -3.51%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.49'
$ws.Range("E38").Value = '  +2.98%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.71'
$ws.Range("E39").Value = '  +1.45%  '

# Row 40
$ws.Range("D40").Value = '0.0₃0753'
$ws.Range("E40").Value = '  +2.44%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.374'
$ws.Range("E41").Value = '  +2.26%  '

# Row 42
$ws.Range("D42").Value = '3.091.52'
$ws.Range("E42").Value = '  -2.32%  '

# Row 43
$ws.Range("E43").Value = '  +0.12%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.84'
$ws.Range("E44").Value = '  +0.96%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0415'
$ws.Range("E45").Value = '  +3.63%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.49'
$ws.Range("E46").Value = '  +2.16%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.21'
$ws.Range("E47").Value = '  +0.11%  '

# Row 48
$ws.Range("E48").Value = '  +5.63%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.58'
$ws.Range("E49").Value = '  -0.67%  '

# Row 50
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.34'
$ws.Range("E50").Value = '  +3.05%  '

# Row 51
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.34'
$ws.Range("E51").Value = '  +4.27%  '
